$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = 'http://www.emelevadores.com.br/'
$ws.Cells.Item(2,2).Value = 'contato@emelevadores.com.br;'
$ws.Cells.Item(2,3).Value = '(11) 2035-1975;'

# Row 3
$ws.Cells.Item(3,1).Value = 'https://www.elevadoreskorman.com.br/empresas-elevadores-sao-paulo'
$ws.Cells.Item(3,2).Value = 'korman@elevadoreskorman.com.br;comercial@elevadoreskorman.com.br;vendas@elevadoreskorman.com.br;'
$ws.Cells.Item(3,3).Value = '(11) 2914-0248;(11) 2914-7506;'

# Row 4
$ws.Cells.Item(4,1).Value = 'https://retrofitelevadores.com.br/'
$ws.Cells.Item(4,2).Value = 'contato@elevadoresretrofit.com.br;'
$ws.Cells.Item(4,3).Value = '(11) 98942-8956;'

# Row 5
$ws.Cells.Item(5,1).Value = 'http://primac.com.br/'
$ws.Cells.Item(5,2).Value = 'comercial@primac.com.br;'
$ws.Cells.Item(5,3).Value = '(11) 2942-7479;'

# Row 6
$ws.Cells.Item(6,1).Value = 'https://iesab.com.br/preco-do-elevador-residencial/#:~:text=M%C3%A9dia%20de%20pre%C3%A7o%20do%20Elevador,comprimento%2C%20menor%20ser%C3%A1%20o%20custo.'
$ws.Cells.Item(6,2).Value = 'emailbit21@gmail.com;'
$ws.Cells.Item(6,3).Value = '(31) 3212-1604;'

# Row 7
$ws.Cells.Item(7,1).Value = 'https://spelevadores.com.br/'
$ws.Cells.Item(7,2).Value = 'contato@spelevadores.com.br;'
$ws.Cells.Item(7,3).Value = '(11) 2353-5320;'

# Row 8
$ws.Cells.Item(8,1).Value = 'https://villarta.com.br/elevadores-e-escadas-rolantes-villarta/lista-de-empresas-de-elevadores-em-sp/'
$ws.Cells.Item(8,2).Value = 'protecaodedados@villarta.com.br;'
$ws.Cells.Item(8,3).Value = '(11) 91364-5830;(11) 3346-8811;'

# Formatting: Arial 10, not bold, left-aligned - matches the new style used by the data rows
$dataRange = $ws.Range("A2:C8")
$dataRange.Font.Name = "Arial"
$dataRange.Font.Size = 10
$dataRange.HorizontalAlignment = -4131
